$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (moves from 2024-04-24 to 2024-05-24, serial 45406 -> 45436)
$ws.Range("A1").Value = "5/24/2024"

# Update the price list values in column D
$ws.Range("D25").Value = 1252.84
$ws.Range("D26").Value = 1396.94
$ws.Range("D27").Value = 1574.396
$ws.Range("D28").Value = 1723.831
$ws.Range("D29").Value = 1894.6
$ws.Range("D30").Value = 2068.056
